$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.692.71'
$ws.Range("E2").Value = '  +3.36%  '
$ws.Range("D3").Value = '3.462.42'
$ws.Range("E3").Value = '  +4.12%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.39'
$ws.Range("E5").Value = '  +4.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '157.05'
$ws.Range("E6").Value = '  +3.66%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.464.02'
$ws.Range("E8").Value = '  +4.09%  '
$ws.Range("E9").Value = '  +7.08%  '
$ws.Range("E10").Value = '  +0.97%  '
$ws.Range("E11").Value = '  +6.66%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  +2.95%  '
$ws.Range("D13").Value = '4.060.36'
$ws.Range("E13").Value = '  +4.23%  '
$ws.Range("E14").Value = '  -1.64%  '
$ws.Range("E15").Value = '  +9.70%  '
$ws.Range("E16").Value = '  +3.73%  '
$ws.Range("D17").Value = '64.700.17'
$ws.Range("E17").Value = '  +3.48%  '
$ws.Range("D18").Value = '3.464.37'
$ws.Range("E18").Value = '  +4.30%  '
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.42'
$ws.Range("E20").Value = '  +4.43%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.58'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '397.71'
$ws.Range("E22").Value = '  +3.54%  '
$ws.Range("E23").Value = '  +1.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.23'
$ws.Range("E24").Value = '  +3.46%  '
$ws.Range("E25").Value = '  -0.48%  '
$ws.Range("E26").Value = '  +25.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.62'
$ws.Range("E27").Value = '  +9.07%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.03'
$ws.Range("E30").Value = '  +8.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.76'
$ws.Range("E31").Value = '  +4.91%  '
$ws.Range("E32").Value = '  +3.34%  '
$ws.Range("E33").Value = '  +5.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.88'
$ws.Range("E34").Value = '  +4.00%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.02'
$ws.Range("E36").Value = '  +4.33%  '
$ws.Range("E37").Value = '  -1.21%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '160.62'
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0786'
$ws.Range("E39").Value = '  +7.03%  '
$ws.Range("E40").Value = '  +0.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '27.60'
$ws.Range("E41").Value = '  +2.83%  '
$ws.Range("D42").Value = '2.909.30'
$ws.Range("E42").Value = '  +2.34%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0324'
$ws.Range("E43").Value = '  +3.32%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.45'
$ws.Range("E44").Value = '  +2.59%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.31'
$ws.Range("E45").Value = '  +4.22%  '
$ws.Range("E46").Value = '  +3.03%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.84'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.09'
$ws.Range("E48").Value = '  +5.58%  '
$ws.Range("E49").Value = '  +24.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.59'
$ws.Range("E50").Value = '  +4.64%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.107'
$ws.Range("E51").Value = '  +2.37%  '

Write-Output "Applied 87 changes"
